# Pathogen.env.1.0 header re-shuffle:
#   - Existing columns G..N (locus_tag_prefix, strain, isolate, collected_by,
#     collection_date, geo_loc_name, isolation_source, lat_lon) shift so that
#     "locus_tag_prefix" moves out to column P, and everything else shifts
#     one column to the left (G..M), making room for the new
#     "culture_collection" column at N.
#   - Ten brand-new attribute columns are appended: culture_collection (N),
#     genotype (O), locus_tag_prefix (P, moved), passage_history (Q),
#     pathotype (R), serotype (S), serovar (T), specimen_voucher (U),
#     subgroup (V), subtype (W).
#   - Each header cell keeps/gets the correct mandatory(green)/
#     at-least-one(blue)/optional(yellow) fill, and the matching cell
#     comment (field definition tooltip).
#
# NOTE: this runtime's PowerShell parser only understands POSITIONAL
# arguments when calling functions (named "-Param value" args come through
# empty), so every helper call below is positional.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cells whose existing format (fill/font/border) represents each
# of the three header classes already present in the sheet.
$greenSample  = $ws.Range("A15")   # mandatory (GREEN)
$blueSample   = $ws.Range("H15")   # at least one mandatory (BLUE)
$yellowSample = $ws.Range("C15")   # optional (YELLOW)

function Set-HeaderCell {
    param([string]$Ref, [string]$HeaderText, [string]$Class, [string]$CommentText)

    $cell = $ws.Range($Ref)
    $cell.Value = $HeaderText

    if ($Class -eq "green") {
        $srcFormat = $greenSample
    } elseif ($Class -eq "blue") {
        $srcFormat = $blueSample
    } else {
        $srcFormat = $yellowSample
    }
    $srcFormat.Copy()
    $cell.PasteSpecial(-4122)   # xlPasteFormats

    $cell.AddComment($CommentText)
}

# --- columns that shift left (old G..N content moves to new G..M, N) -------
Set-HeaderCell "G15" "strain" "blue" "Organism group

microbial or eukaryotic strain name"

Set-HeaderCell "H15" "isolate" "blue" "Organism group

Identification or description of the specific individual from which this sample was obtained"

Set-HeaderCell "I15" "collected_by" "green" "Name of persons or institute who collected the sample"

Set-HeaderCell "J15" "collection_date" "green" "Time of sampling (single instance or interval, eg., 2008-01-23T19:23:10, 2008-01-23, 2008-01, 2008, 1952-10-21T11:43Z/1952-10-21T17:43Z, 1952-10-21/1953-02-15, 1952-10/1953-02, 1952/1953)"

Set-HeaderCell "K15" "geo_loc_name" "green" 'Geographical origin of the sample; use the appropriate name from the list, http://www.ddbj.nig.ac.jp/sub/country-e.html. Use a colon to separate the country or ocean from more detailed information about the location, eg "Japan:Kanagawa, Hakone, Lake Ashi" '

Set-HeaderCell "L15" "isolation_source" "green" "Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived."

Set-HeaderCell "M15" "lat_lon" "green" 'The geographical coordinates of the location where the sample was collected. Specify as decimal degrees latitude and longitude in format "d[d.dddd] N|S d[dd.dddd] W|E", eg, 47.94 N 28.12 W'

# --- new / moved columns ----------------------------------------------------
Set-HeaderCell "N15" "culture_collection" "yellow" "Name of source institute and unique culture identifier. See the description for the proper format and list of allowed institutes, http://www.insdc.org/controlled-vocabulary-culturecollection-qualifier"

Set-HeaderCell "O15" "genotype" "yellow" "observed genotype"

Set-HeaderCell "P15" "locus_tag_prefix" "yellow" "A locus tag prefix required for an annotated genome, http://www.ddbj.nig.ac.jp/sub/locus_tag-e.html"

Set-HeaderCell "Q15" "passage_history" "yellow" "Number of passages and passage method"

Set-HeaderCell "R15" "pathotype" "yellow" "Some bacterial specific pathotypes (example Eschericia coli - STEC, UPEC)"

Set-HeaderCell "S15" "serotype" "yellow" 'Taxonomy below subspecies; a variety (in bacteria, fungi or virus) usually based on its antigenic properties. Same as serovar and serogroup. e.g. serotype="H1N1" in Influenza A virus CY098518.'

Set-HeaderCell "T15" "serovar" "yellow" "Taxonomy below subspecies; a variety (in bacteria, fungi or virus) usually based on its antigenic properties. Same as serovar and serotype. Sometimes used as species identifier in bacteria with shaky taxonomy, e.g. Leptospira interrogans serovar Hardjo, http://www.ncbi.nlm.nih.gov/Taxonomy/Browser/wwwtax.cgi?mode=Info&id=176&lvl=3&lin=f&srchmode=3&unlock"

Set-HeaderCell "U15" "specimen_voucher" "yellow" 'Identifier for the physical specimen. Use format: "[<institution-code>:[<collection-code>:]]<specimen_id>", eg, "UAM:Mamm:52179". Intended as a reference to the physical specimen that remains after it was analyzed. If the specimen was destroyed in the process of analysis, electronic images (e-vouchers) are an adequate substitute for a physical voucher specimen. Ideally the specimens will be deposited in a curated museum, herbarium, or frozen tissue collection, but often they will remain in a personal or laboratory collection for some time before they are deposited in a curated collection. There are three forms of specimen_voucher qualifiers. If the text of the qualifier includes one or more colons it is a ''structured voucher''. Structured vouchers include institution-codes (and optional collection-codes) taken from a controlled vocabulary maintained by the INSDC that denotes the museum or herbarium collection where the specimen resides, please visit the INSDC website, http://www.insdc.org/controlled-vocabulary-specimenvoucher-qualifier'

Set-HeaderCell "V15" "subgroup" "yellow" "Taxonomy below subspecies; sometimes used in viruses to denote subgroups taken from a single isolate."

Set-HeaderCell "W15" "subtype" "yellow" "Used as classifier in viruses (e.g. HIV type 1, Group M, Subtype A)."

Write-Host "Header row updated: G15:W15 now span strain..subtype"
